# Auto-generated edit script: updates crypto price/volume table cells
# to match the target snapshot (commit: "Updated cryptos list on Sun Sep  8 09:24:29 UTC 2024 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    # Force the cell to hold a literal text value even when the text looks
    # like a number (e.g. "503.41"), matching the source data which stores
    # all prices/percentages as text. Reset back to the default/Normal
    # style afterwards so no stray number-format style is left behind.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "54.429.52"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.286.45"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "503.41"
$ws.Range("E5").Value = "  +1.85%  "
Set-TextValue "D6" "130.32"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E7").Value = "  -0.26%  "
Set-TextValue "D8" "0.530"
$ws.Range("E8").Value = "  +0.44%  "
Set-TextValue "D9" "0.0960"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +4.79%  "
Set-TextValue "D12" "4.72"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").Value = "2.696.03"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  +6.59%  "
$ws.Range("D15").Value = "54.438.37"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "2.302.58"
$ws.Range("E17").Value = "  +0.31%  "
Set-TextValue "D18" "10.31"
Set-TextValue "D19" "4.16"
$ws.Range("E19").Value = "  +2.93%  "
Set-TextValue "D20" "305.44"
$ws.Range("E20").Value = "  +0.97%  "
Set-TextValue "D21" "6.42"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +0.14%  "
Set-TextValue "D23" "61.95"
$ws.Range("E23").Value = "  -2.79%  "
Set-TextValue "D24" "0.998"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  +1.82%  "
Set-TextValue "D26" "7.38"
$ws.Range("E26").Value = "  +3.53%  "
Set-TextValue "D27" "173.23"
$ws.Range("E27").Value = "  +4.72%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0696"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D30" "5.99"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("E32").Value = "  -0.01%  "
Set-TextValue "D33" "17.85"
$ws.Range("E33").Value = "  +1.54%  "
Set-TextValue "D34" "0.976"
$ws.Range("E34").Value = "  +12.89%  "
Set-TextValue "D35" "0.995"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +1.84%  "
Set-TextValue "D37" "3.78"
$ws.Range("E37").Value = "  +4.26%  "
Set-TextValue "D38" "0.376"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  +1.66%  "
Set-TextValue "D41" "4.90"
$ws.Range("E41").Value = "  +2.17%  "
Set-TextValue "D42" "125.43"
$ws.Range("E42").Value = "  -0.53%  "
Set-TextValue "D43" "0.0497"
$ws.Range("E43").Value = "  +3.71%  "
Set-TextValue "D44" "0.0896"
$ws.Range("E44").Value = "  +0.85%  "
Set-TextValue "D45" "245.08"
$ws.Range("E45").Value = "  +3.54%  "
Set-TextValue "D46" "0.549"
$ws.Range("E46").Value = "  +0.74%  "
Set-TextValue "D47" "0.374"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +0.89%  "
Set-TextValue "D50" "16.52"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("E51").Value = "  -0.12%  "
